# Applies the Seraph_Profits price/profit refresh described in the commit diff.
# For each affected sheet, update the changed currentAveragePrice* / LevePrice* /
# LeveProfit* cells (columns H-N) to the newly scraped values. A few rows gain or
# lose a trailing LeveProfitHQ (N) / LeveProfitNQ (M) cell entirely, matching the
# upstream diff's cell-level add/remove.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 169.84616
$ws.Range("I33").Value = 175.66667
$ws.Range("K33").Value = 175.66667
$ws.Range("M33").Value = 53.33332999999999
$ws.Range("H80").Value = 269.96667
$ws.Range("I80").Value = 252.6875
$ws.Range("J80").Value = 289.7143
$ws.Range("K80").Value = 758.0625
$ws.Range("L80").Value = 869.1428999999999
$ws.Range("M80").Value = 239.9375
$ws.Range("N80").Value = -2865.1429
$ws.Range("H83").Value = 269.96667
$ws.Range("I83").Value = 252.6875
$ws.Range("J83").Value = 289.7143
$ws.Range("K83").Value = 2274.1875
$ws.Range("L83").Value = 2607.4287
$ws.Range("M83").Value = 2717.8125
$ws.Range("N83").Value = -12591.4287
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -1558
$ws.Range("H137").Value = 2752.4285
$ws.Range("I137").Value = 1303.9166
$ws.Range("K137").Value = 3911.7498
$ws.Range("M137").Value = -1361.7498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3877.8333
$ws.Range("I61").Value = 3751.125
$ws.Range("K61").Value = 3751.125
$ws.Range("M61").Value = -3539.125
$ws.Range("H102").Value = 1664.9286
$ws.Range("I102").Value = 1675
$ws.Range("K102").Value = 1675
$ws.Range("M102").Value = -53
$ws.Range("H136").Value = 3877.8333
$ws.Range("I136").Value = 3751.125
$ws.Range("K136").Value = 11253.375
$ws.Range("M136").Value = -8703.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 700.2
$ws.Range("I22").Value = 678
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 678
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -1246
$ws.Range("H86").Value = 1907.5555
$ws.Range("I86").Value = 1413
$ws.Range("K86").Value = 1413
$ws.Range("M86").Value = -290
$ws.Range("H89").Value = 1907.5555
$ws.Range("I89").Value = 1413
$ws.Range("K89").Value = 7065
$ws.Range("M89").Value = -1449
$ws.Range("H94").Value = 817.0741
$ws.Range("I94").Value = 817.0741
$ws.Range("K94").Value = 817.0741
$ws.Range("M94").Value = -366.0741
$ws.Range("H134").Value = 2109.8518
$ws.Range("I134").Value = 1850.7826
$ws.Range("K134").Value = 5552.3478
$ws.Range("M134").Value = -3017.3478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15506.087
$ws.Range("I99").Value = 14280
$ws.Range("J99").Value = 16160
$ws.Range("K99").Value = 14280
$ws.Range("L99").Value = 16160
$ws.Range("M99").Value = -12782
$ws.Range("N99").Value = -19156
$ws.Range("H105").Value = 1901.2222
$ws.Range("I105").Value = 1587.7142
$ws.Range("J105").Value = 2998.5
$ws.Range("K105").Value = 1587.7142
$ws.Range("L105").Value = 2998.5
$ws.Range("M105").Value = 159.2858000000001
$ws.Range("N105").Value = -6492.5
$ws.Range("H122").Value = 2178.2
$ws.Range("I122").Value = 2245.0454
$ws.Range("K122").Value = 6735.1362
$ws.Range("M122").Value = -4285.1362
$ws.Range("H126").Value = 15506.087
$ws.Range("I126").Value = 14280
$ws.Range("J126").Value = 16160
$ws.Range("K126").Value = 42840
$ws.Range("L126").Value = 48480
$ws.Range("M126").Value = -40370
$ws.Range("N126").Value = -53420

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224
$ws.Range("H51").Value = 833.3333
$ws.Range("I51").Value = 200
$ws.Range("J51").Value = 1150
$ws.Range("K51").Value = 600
$ws.Range("L51").Value = 3450
$ws.Range("M51").Value = -140
$ws.Range("N51").Value = -4370
$ws.Range("H120").Value = 15423.077
$ws.Range("I120").Value = 9000
$ws.Range("J120").Value = 15958.333
$ws.Range("K120").Value = 27000
$ws.Range("L120").Value = 47874.999
$ws.Range("M120").Value = -22162
$ws.Range("N120").Value = -57550.999
$ws.Range("H129").Value = 2098.9167
$ws.Range("I129").Value = 693.1667
$ws.Range("K129").Value = 2079.5001
$ws.Range("M129").Value = 2920.4999
$ws.Range("H131").Value = 687.3333
$ws.Range("I131").Value = 687.3333
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2061.9999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2978.0001
$ws.Range("N131").ClearContents()
$ws.Range("H139").Value = 1322.8889
$ws.Range("I139").Value = 925.875
$ws.Range("K139").Value = 2777.625
$ws.Range("M139").Value = 2362.375
$ws.Range("H140").Value = 60005.332
$ws.Range("I140").Value = 60005.332
$ws.Range("K140").Value = 180015.996
$ws.Range("M140").Value = -174835.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1098.8572
$ws.Range("I107").Value = 1213.7778
$ws.Range("K107").Value = 1213.7778
$ws.Range("M107").Value = 706.2221999999999
$ws.Range("H122").Value = 66345.06
$ws.Range("I122").Value = 2972.8
$ws.Range("K122").Value = 8918.400000000001
$ws.Range("M122").Value = -6468.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2672.1667
$ws.Range("I7").Value = 2706.6
$ws.Range("J7").Value = 2500
$ws.Range("K7").Value = 2706.6
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = -2594.6
$ws.Range("N7").Value = -2724
$ws.Range("H14").Value = 13359
$ws.Range("I14").Value = 36
$ws.Range("J14").Value = 40005
$ws.Range("K14").Value = 36
$ws.Range("L14").Value = 40005
$ws.Range("M14").Value = 136
$ws.Range("N14").Value = -40349
$ws.Range("H22").Value = 8442.571
$ws.Range("I22").Value = 1419.8
$ws.Range("K22").Value = 1419.8
$ws.Range("M22").Value = -1124.8
$ws.Range("H27").Value = 8442.571
$ws.Range("I27").Value = 1419.8
$ws.Range("K27").Value = 1419.8
$ws.Range("M27").Value = -1312.8
$ws.Range("H41").Value = 9250
$ws.Range("I41").Value = 9250
$ws.Range("K41").Value = 9250
$ws.Range("M41").Value = -8812
$ws.Range("H46").Value = 2981.5833
$ws.Range("I46").Value = 2041.3334
$ws.Range("J46").Value = 3545.7334
$ws.Range("K46").Value = 2041.3334
$ws.Range("L46").Value = 3545.7334
$ws.Range("M46").Value = -1853.3334
$ws.Range("N46").Value = -3921.7334
$ws.Range("H68").Value = 1999
$ws.Range("I68").Value = 1999
$ws.Range("K68").Value = 1999
$ws.Range("M68").Value = -1250
$ws.Range("H71").Value = 1999
$ws.Range("I71").Value = 1999
$ws.Range("K71").Value = 9995
$ws.Range("M71").Value = -6251
$ws.Range("H126").Value = 2672.1667
$ws.Range("I126").Value = 2706.6
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 8119.799999999999
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -5649.799999999999
$ws.Range("N126").Value = -12440
$ws.Range("H136").Value = 4665.6665
$ws.Range("I136").Value = 3999
$ws.Range("K136").Value = 11997
$ws.Range("M136").Value = -9447

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6450.852
$ws.Range("I62").Value = 4334.909
$ws.Range("J62").Value = 7905.5625
$ws.Range("K62").Value = 4334.909
$ws.Range("L62").Value = 7905.5625
$ws.Range("M62").Value = -3710.909
$ws.Range("N62").Value = -9153.5625
$ws.Range("H65").Value = 6450.852
$ws.Range("I65").Value = 4334.909
$ws.Range("J65").Value = 7905.5625
$ws.Range("K65").Value = 21674.545
$ws.Range("L65").Value = 39527.8125
$ws.Range("M65").Value = -18554.545
$ws.Range("N65").Value = -45767.8125
$ws.Range("H107").Value = 441.44446
$ws.Range("I107").Value = 337.2
$ws.Range("J107").Value = 481.53845
$ws.Range("K107").Value = 1011.6
$ws.Range("L107").Value = 1444.61535
$ws.Range("M107").Value = 908.4000000000001
$ws.Range("N107").Value = -5284.61535
$ws.Range("H130").Value = 71249.25
$ws.Range("J130").Value = 71249.25
$ws.Range("L130").Value = 71249.25
$ws.Range("N130").Value = -81289.25
$ws.Range("H136").Value = 1639.6923
$ws.Range("I136").Value = 1255.8857
$ws.Range("K136").Value = 3767.6571
$ws.Range("M136").Value = -1217.6571
